$d = $word.ActiveDocument

$pairs = @(
    @("13×51=663", "63×34=2142"),
    @("62×13=806", "32×62=1984"),
    @("69×18=1242", "34×35=1190"),
    @("11×60=660", "81×35=2835"),
    @("35×76=2660", "56×55=3080"),
    @("96×60=5760", "30×57=1710"),
    @("31×29=899", "15×41=615"),
    @("59×20=1180", "73×49=3577"),
    @("64×54=3456", "46×95=4370"),
    @("45×53=2385", "14×33=462"),
    @("81×49=3969", "24×73=1752"),
    @("53×66=3498", "69×40=2760"),
    @("79×81=6399", "86×44=3784"),
    @("84×72=6048", "94×19=1786"),
    @("47×48=2256", "26×84=2184"),
    @("82×28=2296", "69×53=3657"),
    @("42×34=1428", "96×27=2592"),
    @("34×37=1258", "83×36=2988"),
    @("42×14=588", "12×54=648"),
    @("25×93=2325", "36×11=396"),
    @("42×72=3024", "96×40=3840"),
    @("36×60=2160", "79×73=5767"),
    @("85×95=8075", "45×79=3555"),
    @("94×43=4042", "39×83=3237"),
    @("56×24=1344", "29×76=2204")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
